$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a few existing values in the July block (before the new row) ---
$ws.Range("B15").Value2 = 7189.849999999999
$ws.Range("B19").Value2 = 4856.5
$ws.Range("B20").Value2 = 10685.82

# --- Insert a new row at 21 (shifts rows 21:82 down to 22:83) and fill it in ---
$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value2 = 28
$ws.Range("B21").Value2 = 7408.82
$ws.Range("C21").Value2 = 7
$ws.Range("D21").Value2 = 2025
$ws.Range("E21").Value2 = "07/2025"
